$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: header row gains a "Text" label column and an "label" column ---
$ws.Range("C6").Value = 'Above is a well insulated piston cylinder assembly.  Place the block on top of the piston and observe the response.  Calculate the amount of work that the piston and block did on the system?'
$ws.Range("C7").Value = 'The system had an initial temperature of 200 K and contains 1.8 moles of an ideal monatomic gas.  You wrote that get(''WorkDoneAnswer'',''int'') kJ of work were done.  What final temperature should the system have?'
$ws.Range("C8").Value = 'Previously you answered that the compression did get(''WorkDoneAnswer'',''int'') KJ on the system bringing it to a final temperature of get(''TempAnswer'',''int'') K.  Here''s the same compression, but this time we''re displaying work done and temperature. How do the results compare?  If there''s a discrepency, can you account for it?'
$ws.Range("C11").Value = 'Given our ## P_{ext} ## should the graph be linear or did something go wrong? Explain.'
$ws.Range("C2").Value = 'Text'
$ws.Range("E2").Value = 'label'
$ws.Range("C14").Value = 'The system has undergone a two-step process.  First it was compressed by adding a block.  Then it was expanded to its original volume by removing the block.  Before the compression, the system''s temperature was 200K.  After the expansion, the temperature was get(''Temp'', ''int'') K.  Why is the system temperature higher after going through this two-step process?'
$ws.Range("E9").Value = 'Slope from graph'
$ws.Range("E10").Value = 'Slope from equation'
$ws.Range("E12").Value = 'Work Done:'
$ws.Range("E13").Value = 'Final Temperature:'

# --- New answer-choice table starting at row 17 ---
$ws.Range("B17").Value = 'answerId'
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = '## W = -\\int_{V_{1}}^{V_{2}}P_{sys}dV ##'
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = '## W = - V\\Delta P_{ext} ##'
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = '## W = -P_{ext}\\Delta V ##'
$ws.Range("B21").Value = 3
$ws.Range("C21").Value = '## W = -T\\Delta V ##'
$ws.Range("B22").Value = 4
$ws.Range("C22").Value = '## nc_v\\Delta T = Q ##'
$ws.Range("B23").Value = 5
$ws.Range("C23").Value = '##nc_v\\Delta T = -P_{ext}\\Delta V ##'
$ws.Range("B24").Value = 6
$ws.Range("C24").Value = '##nc_p\\Delta T = -P_{ext}\\Delta V ##'
$ws.Range("B25").Value = 7
$ws.Range("C25").Value = 'None of these are correct'

# --- Re-apply the existing "highlight" cell style (same as C5) to the new answer/explanation cells ---
$ws.Range("C8").NumberFormat = "General"
$ws.Range("C9").NumberFormat = "General"
$ws.Range("C10").NumberFormat = "General"
$ws.Range("C14").NumberFormat = "General"

# --- Update the saved view/selection state ---
$ws.Range("G16").Select()
